$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $count = $range.InlineShapes.Count
    if ($count -gt 0) {
        $inlineShape = $range.InlineShapes.Item(1)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape() | Out-Null
    }
}

# Footer (default) and Footer (first page) both hold the Pearson Edexcel
# logo image, currently named "image1.png" -> rename to "image2.png".
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# Header (first page) holds the BTEC logo image, currently named
# "image2.jpg" -> rename to "image1.jpg".
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"
